$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J2: REPORT_TYPE_CODE 001 -> 002 (must remain text, not become numeric 2)
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "002"
$ws.Range("J2").Style = "Normal"

# N2: REPORT_DATE text update
$ws.Range("N2").Value = "2020-06-30 00:00:00"

# Updated figures
$ws.Range("O2").Value = 1260434375.12
$ws.Range("P2").Value = 307629087.06
$ws.Range("Q2").Value = 161524508.23
$ws.Range("R2").Value = ""
$ws.Range("S2").Value = 308737190.24
$ws.Range("T2").Value = ""
$ws.Range("U2").Value = 148196292.68
$ws.Range("V2").Value = ""
$ws.Range("W2").Value = 469545000.39
$ws.Range("X2").Value = 116105376.4
$ws.Range("Y2").Value = ""
$ws.Range("Z2").Value = ""
$ws.Range("AA2").Value = ""
$ws.Range("AB2").Value = 790889374.73
$ws.Range("AC2").Value = ""
$ws.Range("AD2").Value = ""
$ws.Range("AE2").Value = ""
$ws.Range("AF2").Value = 155.4964857511
$ws.Range("AG2").Value = 37.2526336681
